$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dTextCells = @("D5","D6","D8","D9","D11","D12","D14","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D38","D39","D41","D42","D43","D44","D45","D47","D49")
foreach ($addr in $dTextCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '25.951.06'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '1.643.40'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '215.60'
$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("D6").Value = '0.5089'
$ws.Range("E6").Value = '  +0.97%  '

$ws.Range("E7").Value = '  +0.25%  '

$ws.Range("D8").Value = '0.2563'
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = '0.06384'
$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("E10").Value = '  -0.81%  '

$ws.Range("D11").Value = '0.07783'
$ws.Range("E11").Value = '  +0.67%  '

$ws.Range("D12").Value = '4.291'
$ws.Range("E12").Value = '  +0.99%  '

$ws.Range("D13").Value = '1.644.46'
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("D14").Value = '0.5469'
$ws.Range("E14").Value = '  +0.21%  '

$ws.Range("D15").Value = '0.0₅7842'
$ws.Range("E15").Value = '  -0.52%  '

$ws.Range("D16").Value = '64.55'
$ws.Range("E16").Value = '  +0.64%  '

$ws.Range("D17").Value = '26.023.71'
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  +0.17%  '

$ws.Range("D19").Value = '198.05'
$ws.Range("E19").Value = '  -2.52%  '

$ws.Range("D20").Value = '4.467'
$ws.Range("E20").Value = '  +2.14%  '

$ws.Range("D21").Value = '9.977'
$ws.Range("E21").Value = '  +1.01%  '

$ws.Range("D22").Value = '6.055'
$ws.Range("E22").Value = '  +1.46%  '

$ws.Range("D23").Value = '1.006'
$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").Value = '1.878'
$ws.Range("E24").Value = '  -2.07%  '

$ws.Range("D25").Value = '141.50'
$ws.Range("E25").Value = '  +0.40%  '

$ws.Range("D26").Value = '0.1162'
$ws.Range("E26").Value = '  +2.49%  '

$ws.Range("D27").Value = '6.894'
$ws.Range("E27").Value = '  +2.09%  '

$ws.Range("D28").Value = '15.75'
$ws.Range("E28").Value = '  +0.46%  '

$ws.Range("D29").Value = '1.241'
$ws.Range("E29").Value = '  -0.17%  '

$ws.Range("D30").Value = '0.05015'
$ws.Range("E30").Value = '  +1.09%  '

$ws.Range("D31").Value = '3.252'
$ws.Range("E31").Value = '  -0.57%  '

$ws.Range("D32").Value = '3.195'
$ws.Range("E32").Value = '  +0.28%  '

$ws.Range("D33").Value = '1.544'
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("D34").Value = '2.364'
$ws.Range("E34").Value = '  -0.41%  '

$ws.Range("D35").Value = '0.8996'
$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("D36").Value = '2.583'
$ws.Range("E36").Value = '  -1.63%  '

$ws.Range("D37").Value = '1.133.11'
$ws.Range("E37").Value = '  -1.82%  '

$ws.Range("D38").Value = '0.5485'
$ws.Range("E38").Value = '  -2.14%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01557'
$ws.Range("E39").Value = '  -0.64%  '

$ws.Range("B40").Value = 'BabyDogeCoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D40").Value = '0.0₈131'
$ws.Range("E40").Value = '  +12.28%  '

$ws.Range("D41").Value = '1.006'
$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("D42").Value = '2.547'
$ws.Range("E42").Value = '  -0.86%  '

$ws.Range("D43").Value = '5.621'
$ws.Range("E43").Value = '  -0.47%  '

$ws.Range("D44").Value = '0.8165'
$ws.Range("E44").Value = '  +1.21%  '

$ws.Range("D45").Value = '100.22'
$ws.Range("E45").Value = '  +0.41%  '

$ws.Range("D46").Value = '1.778.12'
$ws.Range("E46").Value = '  +0.12%  '

$ws.Range("D47").Value = '0.4541'
$ws.Range("E47").Value = '  +0.24%  '

$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").Value = '54.87'
$ws.Range("E49").Value = '  +0.21%  '

$ws.Range("E50").Value = '  +0.33%  '

$ws.Range("E51").Value = '  +0.45%  '

foreach ($addr in $dTextCells) { $ws.Range($addr).Style = "Normal" }
